$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row that only held the instructor name (no label in column A) is removed,
# shifting every following row up by one.
$ws.Rows.Item(13).Delete()

# After the shift, several cells were left with the wrong (leftover) text and
# need to be corrected to the real syllabus content.
$ws.Range("B10").Value = "787307 - Luis Fernando Figueiredo Faria"
$ws.Range("C10").Value = "787307 - Luis Fernando Figueiredo Faria"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Re-use the already existing "01/01/2018" text (copy it in) so it stays a
# plain text value instead of being re-interpreted as a date.
$ws.Range("B8").Copy($ws.Range("B15"))
$ws.Range("C8").Copy($ws.Range("C15"))

$ws.Range("B18").Value = "787307 - Luis Fernando Figueiredo Faria"
$ws.Range("C18").Value = "787307 - Luis Fernando Figueiredo Faria"

$ws.Range("B19").Value = "-Provas escritas; -participação e conteúdo de trabalho e seminário;"
$ws.Range("C19").Value = "-Provas escritas; -participação e conteúdo de trabalho e seminário;"

$ws.Range("B20").Value = "Média Final = (Prova1 + Prova2 + Nota de Trabalho) / 3Média final mínima de aprovação = 5,0"
$ws.Range("C20").Value = "Média Final = (Prova1 + Prova2 + Nota de Trabalho) / 3Média final mínima de aprovação = 5,0"

$ws.Range("B21").Value = "(Prova escrita + Média Final)/2         Nota Final mínima para aprovação= 5,0"
$ws.Range("C21").Value = "(Prova escrita + Média Final)/2         Nota Final mínima para aprovação= 5,0"
